# Apply scheduled-runner price/profit recalculation updates across all sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 532
$ws.Range("I2").Value = 284.4
$ws.Range("J2").Value = 1151
$ws.Range("K2").Value = 284.4
$ws.Range("L2").Value = 1151
$ws.Range("M2").Value = -171.4
$ws.Range("N2").Value = -1377
# row 6
$ws.Range("H6").Value = 125289.5
$ws.Range("I6").Value = 166719.33
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 500157.99
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -500045.99
$ws.Range("N6").Value = -3224
# row 33
$ws.Range("H33").Value = 208
$ws.Range("I33").Value = 264
$ws.Range("J33").Value = 96
$ws.Range("K33").Value = 264
$ws.Range("L33").Value = 96
$ws.Range("M33").Value = -35
$ws.Range("N33").Value = -554
# row 40
$ws.Range("H40").Value = 2815.1428
$ws.Range("I40").Value = 2725.5
$ws.Range("J40").Value = 2934.6667
$ws.Range("K40").Value = 2725.5
$ws.Range("L40").Value = 2934.6667
$ws.Range("M40").Value = -2550.5
$ws.Range("N40").Value = -3284.6667
# row 43
$ws.Range("H43").Value = 1159.6
$ws.Range("I43").Value = 1079.5
$ws.Range("J43").Value = 1319.8
$ws.Range("K43").Value = 1079.5
$ws.Range("L43").Value = 1319.8
$ws.Range("M43").Value = -1010.5
$ws.Range("N43").Value = -1457.8
# row 69
$ws.Range("H69").Value = 12433527
$ws.Range("J69").Value = 13563257
$ws.Range("L69").Value = 40689771
$ws.Range("N69").Value = -40691519
# row 72
$ws.Range("H72").Value = 12433527
$ws.Range("J72").Value = 13563257
$ws.Range("L72").Value = 122069313
$ws.Range("N72").Value = -122078049
# row 132
$ws.Range("H132").Value = 1837.6097
$ws.Range("I132").Value = 1701.0769
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 5103.2307
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -2573.2307
$ws.Range("N132").Value = -18560
# row 137
$ws.Range("H137").Value = 10287.77
$ws.Range("I137").Value = 9895.200000000001
$ws.Range("K137").Value = 29685.6
$ws.Range("M137").Value = -27135.6
# row 138
$ws.Range("H138").Value = 1256.75
$ws.Range("I138").Value = 1256.75
$ws.Range("K138").Value = 3770.25
$ws.Range("M138").Value = 1369.75
# row 141
$ws.Range("H141").Value = 4908.591
$ws.Range("I141").Value = 5536.316
$ws.Range("J141").Value = 933
$ws.Range("K141").Value = 16608.948
$ws.Range("L141").Value = 2799
$ws.Range("M141").Value = -11428.948
$ws.Range("N141").Value = -13159

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 2108.1082
$ws.Range("I32").Value = 2108.1082
$ws.Range("K32").Value = 2108.1082
$ws.Range("M32").Value = -1821.1082
# row 45
$ws.Range("H45").Value = 1453
$ws.Range("I45").Value = 1270.6666
$ws.Range("K45").Value = 1270.6666
$ws.Range("M45").Value = -893.6666
# row 74
$ws.Range("H74").Value = 1149.8
$ws.Range("I74").Value = 1155.3334
$ws.Range("K74").Value = 1155.3334
$ws.Range("M74").Value = -281.3334
# row 77
$ws.Range("H77").Value = 1149.8
$ws.Range("I77").Value = 1155.3334
$ws.Range("K77").Value = 5776.666999999999
$ws.Range("M77").Value = -1408.666999999999
# row 122
$ws.Range("H122").Value = 1254.2727
$ws.Range("I122").Value = 1254.2727
$ws.Range("K122").Value = 3762.8181
$ws.Range("M122").Value = -1312.8181
# row 132
$ws.Range("H132").Value = 3828.3225
$ws.Range("I132").Value = 3991.4443
$ws.Range("K132").Value = 11974.3329
$ws.Range("M132").Value = -9444.332900000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 138
$ws.Range("H138").Value = 133497.5
$ws.Range("J138").Value = 133197
$ws.Range("L138").Value = 133197
$ws.Range("N138").Value = -143477

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1837.2222
$ws.Range("I31").Value = 1798.3529
$ws.Range("K31").Value = 1798.3529
$ws.Range("M31").Value = -1503.3529
# row 34
$ws.Range("H34").Value = 1837.2222
$ws.Range("I34").Value = 1798.3529
$ws.Range("K34").Value = 1798.3529
$ws.Range("M34").Value = -1596.3529
# row 62
$ws.Range("H62").Value = 7999.3335
$ws.Range("I62").Value = 7999.5
$ws.Range("J62").Value = 7999
$ws.Range("K62").Value = 7999.5
$ws.Range("L62").Value = 7999
$ws.Range("M62").Value = -7375.5
$ws.Range("N62").Value = -9247
# row 65
$ws.Range("H65").Value = 7999.3335
$ws.Range("I65").Value = 7999.5
$ws.Range("J65").Value = 7999
$ws.Range("K65").Value = 39997.5
$ws.Range("L65").Value = 39995
$ws.Range("M65").Value = -36877.5
$ws.Range("N65").Value = -46235
# row 94
$ws.Range("H94").Value = 1637.25
$ws.Range("I94").Value = 1749.6666
$ws.Range("K94").Value = 1749.6666
$ws.Range("M94").Value = -1298.6666
# row 122
$ws.Range("H122").Value = 3444.9375
$ws.Range("I122").Value = 3430.926
$ws.Range("J122").Value = 3520.6
$ws.Range("K122").Value = 10292.778
$ws.Range("L122").Value = 10561.8
$ws.Range("M122").Value = -7842.778
$ws.Range("N122").Value = -15461.8
# row 134
$ws.Range("H134").Value = 2654.087
$ws.Range("I134").Value = 2654.087
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7962.261
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5427.261
$ws.Range("N134").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 113
$ws.Range("H113").Value = 1060.8529
$ws.Range("I113").Value = 577.86957
$ws.Range("K113").Value = 1733.60871
$ws.Range("M113").Value = 436.39129

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 3126.0667
$ws.Range("I102").Value = 3360.0833
$ws.Range("J102").Value = 2190
$ws.Range("K102").Value = 3360.0833
$ws.Range("L102").Value = 2190
$ws.Range("M102").Value = -1738.0833
$ws.Range("N102").Value = -5434

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 94
$ws.Range("H94").Value = 59666.668
$ws.Range("J94").Value = 59666.668
$ws.Range("L94").Value = 59666.668
$ws.Range("N94").Value = -61018.668
# row 122
$ws.Range("H122").Value = 2798.75
$ws.Range("I122").Value = 2648.8333
$ws.Range("K122").Value = 7946.499899999999
$ws.Range("M122").Value = -5496.499899999999
# row 132
$ws.Range("H132").Value = 2169.6128
$ws.Range("I132").Value = 2176.7827
$ws.Range("J132").Value = 2149
$ws.Range("K132").Value = 6530.348100000001
$ws.Range("L132").Value = 6447
$ws.Range("M132").Value = -4000.348100000001
$ws.Range("N132").Value = -11507

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 14
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 4
$ws.Range("M14").Value = 164
# row 45
$ws.Range("H45").Value = 41410.816
$ws.Range("J45").Value = 49410.832
$ws.Range("L45").Value = 49410.832
$ws.Range("N45").Value = -50392.832
# row 107
$ws.Range("H107").Value = 721.6
$ws.Range("I107").Value = 804.7143
$ws.Range("K107").Value = 2414.1429
$ws.Range("M107").Value = -494.1428999999998
# row 132
$ws.Range("H132").Value = 10042.647
$ws.Range("I132").Value = 8784.733
$ws.Range("K132").Value = 26354.199
$ws.Range("M132").Value = -23824.199
